$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 3).Value = 7007
$ws.Cells.Item(7, 5).Value = 289917074
$ws.Cells.Item(14, 3).Value = 110814
$ws.Cells.Item(14, 5).Value = 253241833
$ws.Cells.Item(37, 3).Value = 23047
$ws.Cells.Item(37, 5).Value = 130182720
$ws.Cells.Item(53, 3).Value = 141676
$ws.Cells.Item(53, 5).Value = 590050995
$ws.Cells.Item(56, 3).Value = 11974
$ws.Cells.Item(56, 5).Value = 187807204
$ws.Cells.Item(63, 3).Value = 14338
$ws.Cells.Item(63, 5).Value = 36137488
$ws.Cells.Item(65, 3).Value = 2010
$ws.Cells.Item(65, 5).Value = 13580554
$ws.Cells.Item(70, 3).Value = 15717
$ws.Cells.Item(70, 5).Value = 24651442
$ws.Cells.Item(74, 3).Value = 938
$ws.Cells.Item(74, 5).Value = 4168041
$ws.Cells.Item(79, 3).Value = 116587
$ws.Cells.Item(79, 5).Value = 447322464
$ws.Cells.Item(81, 3).Value = 17431
$ws.Cells.Item(81, 5).Value = 133553933
$ws.Cells.Item(90, 3).Value = 34348
$ws.Cells.Item(90, 5).Value = 67213232
$ws.Cells.Item(91, 3).Value = 151073
$ws.Cells.Item(91, 5).Value = 481901596
$ws.Cells.Item(92, 3).Value = 408954
$ws.Cells.Item(92, 4).Value = 70903
$ws.Cells.Item(92, 5).Value = 1593050860
$ws.Cells.Item(93, 3).Value = 209443
$ws.Cells.Item(93, 5).Value = 1306950791
$ws.Cells.Item(94, 3).Value = 94125
$ws.Cells.Item(94, 5).Value = 914927107
$ws.Cells.Item(95, 3).Value = 50689
$ws.Cells.Item(95, 5).Value = 928906921
$ws.Cells.Item(96, 5).Value = 786973618
$ws.Cells.Item(97, 3).Value = 2150
$ws.Cells.Item(97, 5).Value = 213846316
$ws.Cells.Item(98, 3).Value = 807
$ws.Cells.Item(98, 5).Value = 117420097
$ws.Cells.Item(104, 3).Value = 135210
$ws.Cells.Item(104, 4).Value = 23286
$ws.Cells.Item(104, 5).Value = 272060324
$ws.Cells.Item(106, 3).Value = 18331
$ws.Cells.Item(106, 5).Value = 41271449
$ws.Cells.Item(107, 3).Value = 6384
$ws.Cells.Item(107, 5).Value = 21931803
$ws.Cells.Item(108, 3).Value = 2827
$ws.Cells.Item(108, 5).Value = 18460239
$ws.Cells.Item(113, 3).Value = 8800
$ws.Cells.Item(113, 5).Value = 12651124
$ws.Cells.Item(114, 5).Value = 9073557
$ws.Cells.Item(115, 3).Value = 11680
$ws.Cells.Item(115, 5).Value = 32886288
$ws.Cells.Item(116, 3).Value = 4548
$ws.Cells.Item(116, 5).Value = 20411989
$ws.Cells.Item(122, 3).Value = 8484
$ws.Cells.Item(122, 5).Value = 12669069
$ws.Cells.Item(131, 3).Value = 75580
$ws.Cells.Item(131, 5).Value = 307202531
$ws.Cells.Item(138, 3).Value = 15
$ws.Cells.Item(138, 5).Value = 626897
$ws.Cells.Item(142, 3).Value = 168966
$ws.Cells.Item(142, 5).Value = 681733793
$ws.Cells.Item(165, 3).Value = 83801
$ws.Cells.Item(165, 4).Value = 17112
$ws.Cells.Item(165, 5).Value = 354961515
$ws.Cells.Item(167, 3).Value = 12217
$ws.Cells.Item(167, 5).Value = 105725619
$ws.Cells.Item(168, 3).Value = 6204
$ws.Cells.Item(168, 5).Value = 100524964
$ws.Cells.Item(174, 3).Value = 226077
$ws.Cells.Item(174, 5).Value = 900500557
$ws.Cells.Item(177, 3).Value = 14706
$ws.Cells.Item(177, 5).Value = 251083376